$d = $word.ActiveDocument

$newTexts = @(
    '[[PERSON_52]] – „bez [[PERSON_52]]“, „o [[PERSON_53]]“'
    '[[PERSON_54]] – „bez [[PERSON_54]]“, „k [[PERSON_55]]“'
    '[[PERSON_56]] – „bez [[PERSON_56]]“, „k [[PERSON_56]]“'
    '[[PERSON_57]] – „bez [[PERSON_57]]“, „k [[PERSON_58]]“'
    '[[PERSON_59]] – „bez [[PERSON_59]]“, „k [[PERSON_60]]“'
    '[[PERSON_61]] – „bez [[PERSON_61]]“, „k [[PERSON_61]]“'
    '[[PERSON_62]] – „bez [[PERSON_62]]“, „k [[PERSON_62]]“'
    '[[PERSON_63]] – „bez [[PERSON_63]]“, „k [[PERSON_63]]“'
    '[[PERSON_64]] – „bez [[PERSON_64]]“, „k [[PERSON_65]]“'
    '[[PERSON_66]] – „bez [[PERSON_66]]“, „k [[PERSON_67]]“'
    '[[PERSON_68]] – „bez [[PERSON_68]]“, „k [[PERSON_68]]“'
    '[[PERSON_69]] – „bez [[PERSON_69]]“, „k [[PERSON_69]]“'
    '[[PERSON_70]] – „bez [[PERSON_70]]“, „k [[PERSON_70]]“'
    '[[PERSON_71]] – „bez [[PERSON_71]]“, „k [[PERSON_71]]“'
    '[[PERSON_72]] – „bez [[PERSON_72]]“, „k [[PERSON_72]]“'
    '[[PERSON_73]] – „bez [[PERSON_73]]“, „k [[PERSON_74]]“'
    '[[PERSON_75]] – „bez [[PERSON_75]]“, „k [[PERSON_75]]“'
    '[[PERSON_76]] – „bez [[PERSON_76]]“, „k [[PERSON_76]]“'
    '[[PERSON_77]] – „bez [[PERSON_77]]“, „k [[PERSON_78]]“'
    '[[PERSON_79]] – „bez [[PERSON_79]]“, „k [[PERSON_80]]“'
    '[[PERSON_81]] – „bez [[PERSON_82]]“, „k [[PERSON_83]]“'
    '[[PERSON_84]] – „bez [[PERSON_84]]“, „k [[PERSON_84]]“'
    '[[PERSON_85]] – „bez [[PERSON_85]]“, „k [[PERSON_85]]“'
    '[[PERSON_86]] – „bez [[PERSON_86]]“, „k [[PERSON_87]]“'
    '[[PERSON_88]] – „bez [[PERSON_88]]“, „k [[PERSON_88]]“'
    '[[PERSON_89]] – „bez [[PERSON_89]]“, „k [[PERSON_89]]“'
    '[[PERSON_90]] – „bez [[PERSON_90]]“, „k [[PERSON_91]]“'
    '[[PERSON_92]] – „bez [[PERSON_92]]“, „k [[PERSON_92]]“'
    '[[PERSON_93]] – „bez [[PERSON_93]]“, „k [[PERSON_93]]“'
    '[[PERSON_94]] – „bez [[PERSON_95]]“, „k [[PERSON_96]]“'
    '[[PERSON_97]] – „bez [[PERSON_97]]“, „k [[PERSON_98]]“'
    '[[PERSON_99]] – „bez [[PERSON_99]]“, „k [[PERSON_99]]“'
    '[[PERSON_100]] – „bez [[PERSON_100]]“, „k [[PERSON_100]]“'
    '[[PERSON_101]] – „bez [[PERSON_102]]“, „k [[PERSON_101]]“'
    '[[PERSON_103]] – „bez [[PERSON_103]]“, „k [[PERSON_103]]“'
    '[[PERSON_104]] – „bez [[PERSON_104]]“, „k [[PERSON_105]]“'
    '[[PERSON_106]] – „bez [[PERSON_106]]“, „k [[PERSON_106]]“'
    '[[PERSON_107]] – „bez [[PERSON_107]]“, „k [[PERSON_107]]“'
    '[[PERSON_108]] – „bez [[PERSON_108]]“, „k [[PERSON_108]]“'
    '[[PERSON_109]] – „bez [[PERSON_109]]“, „k [[PERSON_109]]“'
    '[[PERSON_110]] – „bez [[PERSON_110]]“, „k [[PERSON_110]]“'
    '[[PERSON_111]] – „bez [[PERSON_111]]“, „k [[PERSON_111]]“'
    '[[PERSON_112]] – „bez [[PERSON_112]]“, „k [[PERSON_113]]“'
    '[[PERSON_114]] – „bez [[PERSON_114]]“, „k [[PERSON_114]]“'
    '[[PERSON_115]] – „bez [[PERSON_115]]“, „k [[PERSON_116]]“'
    '[[PERSON_117]] – „bez [[PERSON_117]]“, „k [[PERSON_117]]“'
    '[[PERSON_118]] – „bez [[PERSON_118]]“, „k [[PERSON_118]]“'
    '[[PERSON_119]] – „bez [[PERSON_119]]“, „k [[PERSON_119]]“'
    '[[PERSON_120]] – „bez [[PERSON_120]]“, „k [[PERSON_120]]“'
)

$startPara = 55
$count = 0
for ($i = 0; $i -lt $newTexts.Length; $i++) {
    $paraIndex = $startPara + $i
    $para = $d.Paragraphs($paraIndex)
    $rng = $para.Range
    # Preserve the trailing paragraph mark by only replacing the run text portion
    $endRng = $d.Range($rng.Start, $rng.End - 1)
    $endRng.Text = $newTexts[$i]
    $count = $count + 1
}
Write-Host "Total paragraphs updated:" $count
